$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.735.80'
$ws.Range("E2").Value = '  +3.58%  '

$ws.Range("D3").Value = '2.711.51'
$ws.Range("E3").Value = '  +3.17%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'528.79"
$ws.Range("E5").Value = '  +1.60%  '

$ws.Range("D6").Value = "'147.23"
$ws.Range("E6").Value = '  +1.78%  '

$ws.Range("E7").Value = '  -0.18%  '

$ws.Range("E8").Value = '  +0.46%  '

$ws.Range("D9").Value = '2.727.37'
$ws.Range("E9").Value = '  +3.38%  '

$ws.Range("E10").Value = '  +13.13%  '

$ws.Range("E11").Value = '  +1.52%  '

$ws.Range("E12").Value = '  +2.19%  '

$ws.Range("D14").Value = '3.188.52'
$ws.Range("E14").Value = '  +3.16%  '

$ws.Range("D15").Value = '60.722.12'
$ws.Range("E15").Value = '  +3.52%  '

$ws.Range("D16").Value = "'21.44"
$ws.Range("E16").Value = '  +3.37%  '

$ws.Range("D17").Value = '2.731.77'
$ws.Range("E17").Value = '  +3.72%  '

$ws.Range("E18").Value = '  +1.67%  '

$ws.Range("D19").Value = "'344.69"
$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").Value = "'4.50"
$ws.Range("E20").Value = '  +1.62%  '

$ws.Range("E21").Value = '  +3.61%  '

$ws.Range("D22").Value = "'6.43"
$ws.Range("E22").Value = '  +4.92%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").Value = "'63.36"
$ws.Range("E24").Value = '  +2.69%  '

$ws.Range("E25").Value = '  +4.60%  '

$ws.Range("E26").Value = '  +1.22%  '

$ws.Range("D27").Value = "'0.995"
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D28").Value = '0.0₃0823'
$ws.Range("E28").Value = '  +2.91%  '

$ws.Range("E29").Value = '  +4.26%  '

$ws.Range("D30").Value = "'6.75"
$ws.Range("E30").Value = '  +8.62%  '

$ws.Range("D32").Value = "'1.61"
$ws.Range("E32").Value = '  +2.13%  '

$ws.Range("D33").Value = "'19.07"
$ws.Range("E33").Value = '  +1.34%  '

$ws.Range("D34").Value = "'150.06"
$ws.Range("E34").Value = '  +0.59%  '

$ws.Range("E36").Value = '  +7.85%  '

$ws.Range("D37").Value = "'0.923"
$ws.Range("E37").Value = '  -5.59%  '

$ws.Range("D38").Value = "'0.902"
$ws.Range("E38").Value = '  +7.99%  '

$ws.Range("E39").Value = '  +8.32%  '

$ws.Range("D40").Value = "'37.29"
$ws.Range("E40").Value = '  +1.87%  '

$ws.Range("E41").Value = '  +1.63%  '

$ws.Range("E42").Value = '  +4.35%  '

$ws.Range("D43").Value = "'281.09"
$ws.Range("E43").Value = '  +0.58%  '

$ws.Range("D44").Value = "'20.20"
$ws.Range("E44").Value = '  +3.71%  '

$ws.Range("D45").Value = "'0.995"
$ws.Range("E45").Value = '  -0.24%  '

$ws.Range("D46").Value = "'0.0986"
$ws.Range("E46").Value = '  +0.27%  '

$ws.Range("D47").Value = '2.108.11'
$ws.Range("E47").Value = '  +6.21%  '

$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").Value = "'4.95"
$ws.Range("E48").Value = '  +7.17%  '

$ws.Range("B49").Value = 'Hedera'
$ws.Range("C49").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D49").Value = "'0.0542"
$ws.Range("E49").Value = '  +4.06%  '

$ws.Range("E50").Value = '  +2.19%  '

$ws.Range("D51").Value = "'19.32"
$ws.Range("E51").Value = '  +5.60%  '
